$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.251.16'
$ws.Range('E2').Value = '  +1.83%  '
$ws.Range('D3').Value = '3.540.58'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.38'
$ws.Range('E5').Value = '  +1.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.60'
$ws.Range('E6').Value = '  +1.00%  '
$ws.Range('D7').Value = '3.528.71'
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.195'
$ws.Range('E10').Value = '  +3.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.51'
$ws.Range('E11').Value = '  +10.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.583'
$ws.Range('E12').Value = '  +1.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '46.26'
$ws.Range('E13').Value = '  -2.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000275'
$ws.Range('D15').Value = '4.119.30'
$ws.Range('E15').Value = '  +0.84%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.29'
$ws.Range('E16').Value = '  -1.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '608.42'
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D18').Value = '3.541.88'
$ws.Range('E18').Value = '  +0.70%  '
$ws.Range('D19').Value = '70.314.05'
$ws.Range('E19').Value = '  +1.78%  '
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.25'
$ws.Range('E21').Value = '  -0.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.875'
$ws.Range('E22').Value = '  -1.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.16'
$ws.Range('E23').Value = '  -17.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '15.64'
$ws.Range('E24').Value = '  -0.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '96.32'
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('E26').Value = '  -3.28%  '
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.59'
$ws.Range('E28').Value = '  -0.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.73'
$ws.Range('E29').Value = '  +3.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.99'
$ws.Range('E30').Value = '  -2.14%  '
$ws.Range('B31').Value = 'Stacks'
$ws.Range('C31').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.03'
$ws.Range('E31').Value = '  -2.92%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.20'
$ws.Range('E32').Value = '  -3.78%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '667.27'
$ws.Range('E33').Value = '  +8.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.09'
$ws.Range('E34').Value = '  +2.85%  '
$ws.Range('E35').Value = '  -2.11%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.100'
$ws.Range('E36').Value = '  -1.69%  '
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.57'
$ws.Range('E37').Value = '  +2.86%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '10.72'
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0473'
$ws.Range('E39').Value = '  +6.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '57.04'
$ws.Range('E40').Value = '  +0.10%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E42').Value = '  +4.15%  '
$ws.Range('D43').Value = '3.374.16'
$ws.Range('E43').Value = '  -0.52%  '
$ws.Range('E44').Value = '  -2.64%  '
$ws.Range('D45').Value = '0.0₃0695'
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '32.44'
$ws.Range('E46').Value = '  -0.95%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.92'
$ws.Range('E47').Value = '  +6.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.60'
$ws.Range('E48').Value = '  +2.74%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.129'
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.29'
$ws.Range('E50').Value = '  -1.37%  '
$ws.Range('E51').Value = '  -0.05%  '
